# Labtainer tcpip lab manual: extend the TCP Session Hijacking task (3.4)
# with hints about nping's "-data" usage and what to actually do with the
# hijacked telnet session.
#
# In the original markup the sentence lives in two adjacent <w:r> runs
# that happen to share identical run properties:
#   run A: ") to perform this task.  "
#   run B: "Use the <96>data option ... ascii text."
# The commit appends a trailing space to run A's text (after folding run
# B's original sentence into it) and replaces run B's text with the new
# guidance sentence. We reproduce that by doing the text-level edits with
# Find/Replace, working on the second (later) run first so the search
# string for the first run's replacement doesn't re-match text we just
# inserted.

$d = $word.ActiveDocument

$dash = [char]0x2013   # "-" (en dash), as used by "-data"
$lq   = [char]0x201C   # left curly quote
$rq   = [char]0x201D   # right curly quote

$newSentence = "You will also want to provide the psh and ack flags, and ack the previous packet in your spoofed packet. Your goal is to use a spoofed packet to hijack a telnet session and delete the file on the server at ~/documents/delete-this.txt.  Note that if you use your telnet session to delete that file, e.g., to observe the protocol in wireshark, then you must recreate that file so it can be deleted in a hijacked session."

# Guard against re-running this script against an already-edited document.
$already = $d.Content.Find.Execute($newSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $already) {

    # --- Run B: replace the old "Use the ... ascii text." sentence with
    # the new guidance about psh/ack flags and the hijacking goal. ---
    $oldRunB = "Use the " + $dash + "data option to send your payload.  Your attacker home directory includes a " + $lq + "hexify.py" + $rq + " script that creates hex versions of ascii text."

    $foundB = $d.Content.Find.Execute($oldRunB, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)
    if (-not $foundB) {
        throw "Could not find the original 'Use the -data option...' sentence to replace."
    }

    # --- Run A: keep its own sentence, fold in the (now-superseded)
    # original wording that used to be run B, and add the new trailing
    # space. ---
    $oldRunA = ") to perform this task.  "
    $newRunA = ") to perform this task.  Use the " + $dash + "data option to send your payload.  Your attacker home directory includes a " + $lq + "hexify.py" + $rq + " script that creates hex versions of ascii text. "

    $foundA = $d.Content.Find.Execute($oldRunA, $true, $false, $false, $false, $false, $true, 1, $false, $newRunA, 2)
    if (-not $foundA) {
        throw "Could not find the original ') to perform this task.' text to replace."
    }
}
